# Fix spelling: "veranderd" -> "verandert" in the sentence about
# trajectory length changing per option. This is the only actual
# textual change in the commit; the rest of the underlying XML diff
# is just run-splitting/merging (spell-check proofErr markers removed)
# which Word performs naturally when text is (re)typed/edited.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "verbindingen per traject. Dit veranderd per lengte",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "verbindingen per traject. Dit verandert per lengte",
    2
)
